$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object "object[,]" 1,23
$row2[0,0] = 0.0226437462951986
$row2[0,1] = 0.0224066390041494
$row2[0,2] = 0.0228808535862478
$row2[0,3] = 0.0036751630112626
$row2[0,4] = 0.0002371072910492
$row2[0,5] = 0.0340248962655602
$row2[0,6] = 0.897213989330172
$row2[0,7] = 0.0032009484291642
$row2[0,8] = 0.987077652637819
$row2[0,9] = 0.967753408417309
$row2[0,10] = 0.0469472436277416
$row2[0,11] = 0.000948429164196799
$row2[0,12] = 0.0126852400711322
$row2[0,13] = 0.651570835803201
$row2[0,14] = 0.996324836988737
$row2[0,15] = 0.0002371072910492
$row2[0,16] = 0.0013040901007706
$row2[0,17] = 0.014226437462952
$row2[0,18] = 0
$row2[0,19] = 0.891049199762893
$row2[0,20] = 0.0026081802015412
$row2[0,21] = 0.0020154119739182
$row2[0,22] = 0.157913455838767
$ws.Range("B2:X2").Value = $row2

$row3 = New-Object "object[,]" 1,23
$row3[0,0] = 0.975933609958506
$row3[0,1] = 0.0162418494368702
$row3[0,2] = 0.00580912863070539
$row3[0,3] = 0.972021339656194
$row3[0,4] = 0.950800237107291
$row3[0,5] = 0.794783639596918
$row3[0,6] = 0.0461173681090694
$row3[0,7] = 0.958269116775341
$row3[0,8] = 0.0013040901007706
$row3[0,9] = 0.00663900414937759
$row3[0,10] = 0.0407824540604624
$row3[0,11] = 0.0001185536455246
$row3[0,12] = 0.0001185536455246
$row3[0,13] = 0.0050978067575578
$row3[0,14] = 0.000592768227622999
$row3[0,15] = 0
$row3[0,16] = 0.0002371072910492
$row3[0,17] = 0.984943687018376
$row3[0,18] = 0.999644339063426
$row3[0,19] = 0.000948429164196799
$row3[0,20] = 0.0158861885002964
$row3[0,21] = 0.0224066390041494
$row3[0,22] = 0.000711321873147599
$ws.Range("B3:X3").Value = $row3

$row4 = New-Object "object[,]" 1,23
$row4[0,0] = 0.0003556609365738
$row4[0,1] = 0.0003556609365738
$row4[0,2] = 0.0024896265560166
$row4[0,3] = 0.0002371072910492
$row4[0,4] = 0.0232365145228216
$row4[0,5] = 0.0001185536455246
$row4[0,6] = 0.0425607587433314
$row4[0,7] = 0.0359217545939538
$row4[0,8] = 0.00497925311203319
$row4[0,9] = 0.0212211025489034
$row4[0,10] = 0.9036158861885
$row4[0,11] = 0.998814463544754
$row4[0,12] = 0.986840545346769
$row4[0,13] = 0.323532898636633
$row4[0,14] = 0.000711321873147599
$row4[0,15] = 0.999762892708951
$row4[0,16] = 0.998103141671606
$row4[0,17] = 0.0001185536455246
$row4[0,18] = 0.0001185536455246
$row4[0,19] = 0.10231179608773
$row4[0,20] = 0.97474807350326
$row4[0,21] = 0.975222288085359
$row4[0,22] = 0.829994072317724
$ws.Range("B4:X4").Value = $row4

$row5 = New-Object "object[,]" 1,23
$row5[0,0] = 0.0003556609365738
$row5[0,1] = 0.95850622406639
$row5[0,2] = 0.967753408417309
$row5[0,3] = 0.0235921754593954
$row5[0,4] = 0.0239478363959692
$row5[0,5] = 0.163604030823948
$row5[0,6] = 0.0114997036158862
$row5[0,7] = 0.0014226437462952
$row5[0,8] = 0.00616478956727919
$row5[0,9] = 0.0002371072910492
$row5[0,10] = 0.000829875518672199
$row5[0,11] = 0
$row5[0,12] = 0
$row5[0,13] = 0.0002371072910492
$row5[0,14] = 0
$row5[0,15] = 0
$row5[0,16] = 0
$row5[0,17] = 0.0004742145820984
$row5[0,18] = 0.0002371072910492
$row5[0,19] = 0.0002371072910492
$row5[0,20] = 0.00675755779490219
$row5[0,21] = 0.0002371072910492
$row5[0,22] = 0.0002371072910492
$ws.Range("B5:X5").Value = $row5

